$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (blank) column before column N (14th column), shifting the
# old N:P ("Late" / "heading"+"Original" / "Outstanding") columns one to the
# right, into O:Q.
$mWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $mWidth

# Set the active cell/selection as it ended up after the edit
$ws.Range("M14").Select()
